# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting from an existing header cell (A1) so the
# new header cells pick up the same bold/bordered/centered style used
# by the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (rows 2-60): team record values for every player row ---
$ws.Range("AD2:AD60").Value = 90
$ws.Range("AE2:AE60").Value = 72
$ws.Range("AF2:AF60").Value = 0
